$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Populate column G (Japanese translation) for each data row, mirroring
# the English entry already present in column B on that row.
$ws.Range("G2").Value = "胎膜早期破裂に伴う早産"
$ws.Range("G3").Value = "中～後期早産"
$ws.Range("G4").Value = "超早産"
$ws.Range("G5").Value = "子宮頸管機能不全または膜脆弱性による早産"
$ws.Range("G6").Value = "超早産"
$ws.Range("G7").Value = "妊娠後期"
$ws.Range("G8").Value = "羊水過少症"
$ws.Range("G9").Value = "メコニウム染色羊水"
$ws.Range("G10").Value = "多羊水膜症"
$ws.Range("G11").Value = "無月経"
$ws.Range("G12").Value = "胎児多尿"
$ws.Range("G13").Value = "胞状奇胎"
$ws.Range("G14").Value = "臍帯分離遅延"
$ws.Range("G15").Value = "四枝臍帯"
$ws.Range("G16").Value = "単一臍動脈"
$ws.Range("G17").Value = "臍動脈の拍動性の上昇"
$ws.Range("G18").Value = "拡張末期臍動脈逆流"
$ws.Range("G19").Value = "拡張末期臍動脈血流の消失"
$ws.Range("G20").Value = "臍帯静脈血流の異常"
$ws.Range("G21").Value = "毛状臍帯挿入"
$ws.Range("G22").Value = "限界臍帯挿入"
$ws.Range("G23").Value = "前置静脈瘤"
$ws.Range("G24").Value = "蛇行臍帯挿入"
$ws.Range("G25").Value = "短い臍帯"
$ws.Range("G26").Value = "臍帯結び目"
$ws.Range("G27").Value = "腹外臍静脈瘤"
$ws.Range("G28").Value = "胎児腹腔内臍静脈瘤"
$ws.Range("G29").Value = "長い臍帯"
$ws.Range("G30").Value = "臍帯血低酸素血症"
$ws.Range("G31").Value = "臍帯血高酸素血症"
$ws.Range("G32").Value = "臍帯血高酸素血症"
$ws.Range("G33").Value = "臍帯血低カプニア血症"
$ws.Range("G34").Value = "臍帯嚢胞"
$ws.Range("G35").Value = "臍帯血腫"
$ws.Range("G36").Value = "新生児脳炎"
$ws.Range("G37").Value = "胎盤剥離"
$ws.Range("G38").Value = "水腫性胎盤"
$ws.Range("G39").Value = "胎盤の厚さ増加"
$ws.Range("G40").Value = "小さな胎盤"
$ws.Range("G41").Value = "絨毛膜下隔嚢胞"
$ws.Range("G42").Value = "胎盤梗塞"
$ws.Range("G43").Value = "胎盤間葉系異形成"
$ws.Range("G44").Value = "浮腫性絨毛"
$ws.Range("G45").Value = "石灰化胎盤"
$ws.Range("G46").Value = "胎盤周囲炎"
$ws.Range("G47").Value = "急性胎盤"
$ws.Range("G48").Value = "増多胎盤"
$ws.Range("G49").Value = "母体血管不全"
$ws.Range("G50").Value = "絨毛血管腫"
$ws.Range("G51").Value = "絨毛膜下血栓血腫"
$ws.Range("G52").Value = "胎児血管不全"
$ws.Range("G53").Value = "絨毛がん"
$ws.Range("G54").Value = "自然羊膜剥離"
$ws.Range("G55").Value = "足の羊膜攣縮輪"
$ws.Range("G56").Value = "趾球狭窄輪"
$ws.Range("G57").Value = "腕の羊膜攣縮輪"
$ws.Range("G58").Value = "慢性絨毛膜炎"
$ws.Range("G59").Value = "羊膜シート"
$ws.Range("G60").Value = "着床前致死"
$ws.Range("G61").Value = "接合卵割不全"
$ws.Range("G62").Value = "反復着床障害"
$ws.Range("G63").Value = "受精時の多発前核形成"
$ws.Range("G64").Value = "誘発経膣分娩"
$ws.Range("G65").Value = "完全逆子"
$ws.Range("G66").Value = "不完全逆子"
$ws.Range("G67").Value = "フランク逆子"
$ws.Range("G68").Value = "二次帝王切開"
$ws.Range("G69").Value = "一次帝王切開"
$ws.Range("G70").Value = "オドン式分娩"
$ws.Range("G71").Value = "鉗子分娩"
$ws.Range("G72").Value = "頭血腫"
$ws.Range("G73").Value = "帝王切開後の膣分娩"
$ws.Range("G74").Value = "肩甲難産"
$ws.Range("G75").Value = "治療的人工妊娠中絶"
$ws.Range("G77").Value = "へその緒"
$ws.Range("G78").Value = "子宮外妊娠"
$ws.Range("G79").Value = "母体高フェニルアラニン血症"
$ws.Range("G80").Value = "自然流産の再発"
$ws.Range("G81").Value = "妊娠悪阻"
$ws.Range("G82").Value = "母体高血圧症"
$ws.Range("G83").Value = "子癇"
$ws.Range("G84").Value = "子癇前症"
$ws.Range("G85").Value = "妊娠中の母体の男性化"
$ws.Range("G86").Value = "低母体循環PAPP-A濃度"
$ws.Range("G87").Value = "高母体循環α-フェト蛋白濃度"
$ws.Range("G88").Value = "低母体循環エストリオール濃度"
$ws.Range("G89").Value = "低母体循環絨毛性ゴナドトロピン濃度"
$ws.Range("G90").Value = "母体循環絨毛性ゴナドトロピン濃度高値"
$ws.Range("G91").Value = "低母体循環αフェトプロテイン濃度"
$ws.Range("G92").Value = "母親の催奇形性曝露"
$ws.Range("G93").Value = "母体血栓症"
$ws.Range("G94").Value = "母体糖尿病"
$ws.Range("G95").Value = "母体第一期発熱"
$ws.Range("G96").Value = "分娩中の発熱"
$ws.Range("G97").Value = "歪んだ母体X不活性化"
$ws.Range("G98").Value = "母体抗カルジオリピン抗体陽性"
$ws.Range("G99").Value = "母体の痙攣"
$ws.Range("G100").Value = "10分間APGARスコア2"
$ws.Range("G101").Value = "10分間APGARスコア5"
$ws.Range("G102").Value = "10分間APGARスコア1"
$ws.Range("G103").Value = "10分間APGARスコア0"
$ws.Range("G104").Value = "10分間のAPGARスコア4"
$ws.Range("G105").Value = "10分間のAPGARスコア6"
$ws.Range("G106").Value = "10分間のAPGARスコア3"
$ws.Range("G107").Value = "1分間のAPGARスコア5"
$ws.Range("G108").Value = "1分間のAPGARスコア6"
$ws.Range("G109").Value = "1分間のAPGARスコア0"
$ws.Range("G110").Value = "1分間のAPGARスコア1"
$ws.Range("G111").Value = "1分間のAPGARスコア3"
$ws.Range("G112").Value = "1分間のAPGARスコア4"
$ws.Range("G113").Value = "1分間のAPGARスコア2"
$ws.Range("G114").Value = "5分間のAPGARスコア0"
$ws.Range("G115").Value = "5分間のAPGARスコア5"
$ws.Range("G116").Value = "5分間のAPGARスコア1"
$ws.Range("G117").Value = "5分間のAPGARスコア4"
$ws.Range("G118").Value = "5分間のAPGARスコア6"
$ws.Range("G119").Value = "5分間のAPGARスコア3"
$ws.Range("G120").Value = "5分間のAPGARスコア2"
$ws.Range("G121").Value = "分APGARスコア2"
$ws.Range("G122").Value = "両側胎児無気肺"
$ws.Range("G123").Value = "エコー源性心内フォーカス"
$ws.Range("G124").Value = "短い胎児上腕骨長"
$ws.Range("G125").Value = "肥大した胎児脳槽"
$ws.Range("G126").Value = "胎児腸管エコー"
$ws.Range("G127").Value = "胎児第5指臨床指節症"
$ws.Range("G128").Value = "軽度の胎児脳室肥大"
$ws.Range("G129").Value = "短い胎児大腿骨長"
$ws.Range("G130").Value = "胎児超音波検査における胃の気泡の欠如"
$ws.Range("G131").Value = "胎児鼻骨の低形成"
$ws.Range("G132").Value = "胎児脈絡叢嚢胞"
$ws.Range("G133").Value = "胎児頚部浮腫"
$ws.Range("G134").Value = "胎児頚部腫瘤"
$ws.Range("G135").Value = "頸部リンパ嚢腫"
$ws.Range("G136").Value = "頚部透光性の増大"
$ws.Range("G137").Value = "胎児嚢胞性水腫"
$ws.Range("G138").Value = "低形成鼻骨"
$ws.Range("G139").Value = "胎児鼻骨の欠如"
$ws.Range("G140").Value = "胎児大腿骨/足長比の減少"
$ws.Range("G141").Value = "シャンパンコルク徴候"
$ws.Range("G142").Value = "レモン徴候"
$ws.Range("G143").Value = "胎児三頭筋"
$ws.Range("G144").Value = "胎児長骨エコー原性の減少"
$ws.Range("G145").Value = "胎児長骨エコー増加"
$ws.Range("G146").Value = "胎児脳実質出血"
$ws.Range("G147").Value = "脳室肥大を伴わない胎児脳室内出血"
$ws.Range("G148").Value = "脳室肥大を伴う胎児脳室内出血"
$ws.Range("G149").Value = "胎児性脳下垂体出血"
$ws.Range("G150").Value = "脳室周囲出血を伴う胎児脳室内出血"
$ws.Range("G151").Value = "胎児軸外出血"
$ws.Range("G152").Value = "胎児小脳半球出血"
$ws.Range("G153").Value = "胎児小脳縁出血"
$ws.Range("G154").Value = "透明腔欠損"
$ws.Range("G155").Value = "胎児水晶体エコー増加"
$ws.Range("G156").Value = "胎児心臓流出路異常"
$ws.Range("G157").Value = "胎児心臓の異常な4室像"
$ws.Range("G158").Value = "卵円孔瘤"
$ws.Range("G159").Value = "胎児心嚢液貯留"
$ws.Range("G160").Value = "静脈管形成不全"
$ws.Range("G161").Value = "胎児頭皮腫瘤"
$ws.Range("G162").Value = "胎児皮膚浮腫"
$ws.Range("G163").Value = "出生直後の皮膚の落屑"
$ws.Range("G164").Value = "カゼ性胎脂様落屑"
$ws.Range("G165").Value = "欠落性産毛"
$ws.Range("G166").Value = "2型先天性肺気道奇形"
$ws.Range("G167").Value = "1型先天性肺気道奇形"
$ws.Range("G168").Value = "3型先天性肺気道奇形"
$ws.Range("G169").Value = "胎児気胸"
$ws.Range("G170").Value = "胎児水胸"
$ws.Range("G171").Value = "胃仮死"
$ws.Range("G172").Value = "胎児腹水"
$ws.Range("G173").Value = "出生前二重気泡徴候"
$ws.Range("G174").Value = "胎児腹部嚢胞"
$ws.Range("G175").Value = "胎児腸管拡張"
$ws.Range("G176").Value = "メコニウム仮性嚢胞"
$ws.Range("G177").Value = "胎児下部尿路閉塞"
$ws.Range("G178").Value = "超音波非観血的胎児膀胱"
$ws.Range("G179").Value = "女性胎児の処女化"
$ws.Range("G180").Value = "非免疫性胎児水腫"
$ws.Range("G181").Value = "胎動の増加"
$ws.Range("G182").Value = "胎児アキネジアシークエンス"
$ws.Range("G183").Value = "胎児貧血"
$ws.Range("G184").Value = "中大脳動脈拍動指数の低下"
$ws.Range("G185").Value = "持続性胎児循環"
$ws.Range("G186").Value = "先天性ポートコステティック静脈シャント"
$ws.Range("G187").Value = "持続性静脈管開存症"
$ws.Range("G188").Value = "静脈管血流異常"
$ws.Range("G189").Value = "胎児中大脳動脈収縮期ピーク速度上昇"
$ws.Range("G190").Value = "妊娠中期以降の片方の双子の子宮内胎児死亡"
$ws.Range("G191").Value = "胎児の苦痛"
$ws.Range("G192").Value = "双胎間輸血"
$ws.Range("G193").Value = "動脈管早期閉鎖"

# Select the full first row (header), matching the saved selection state.
$ws.Range("A1:XFD1").Select()

